$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day of price data (2026-02-15) was published; it is prepended as the
# new row 2, pushing every existing data row down by one. The values for the
# new day repeat the constant series values (783.5 / 1112 / 3610) used
# throughout the whole history.
$ws.Rows.Item(2).Insert()

# Format the date cell as Text first so Excel does not silently convert the
# ISO-looking date string into a date serial number (the source data stores
# dates as plain text).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-15"

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
